$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new columns (BG and BH) with headers and descriptions,
# matching the style (wrap text) of the existing header cells.
$ws.Range("BG1").Value = "late_res"
$ws.Range("BH1").Value = "slow_mvmnt"
$ws.Range("BG2").Value = "subject started moving too late."
$ws.Range("BH2").Value = "movement was too slow"

$ws.Range("BG1:BH2").WrapText = $true

# Scroll the view so column AT is the left-most visible column, and select
# the cell just below the newly added data, matching the end-state view
# recorded in the saved workbook.
$excel.ActiveWindow.ScrollColumn = 46
$ws.Range("BH3").Select()

$wb.Save()
